# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" on every sheet
# - Refresh the associated "Latest HO Xliff Generate Date" / "Latest Handoff
#   Datetime" timestamps to the new handoff run time
# - Column widths for the status columns grow to fit the new (longer) text

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-24 10:40:50"
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet -------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-24 10:40:45"
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet -------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-24 10:40:50"
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
